$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.566.09"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.493.31"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'569.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'165.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "2.491.86"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "2.946.87"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "69.383.19"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'24.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "2.496.52"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "'11.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "'347.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'69.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "'3.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "2.621.44"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'8.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.13%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "0.0₃0871"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "'7.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").Value = "'435.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.55%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("D36").Value = "'154.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").Value = "'0.113"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("D38").Value = "'19.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "'18.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("B44").Value = "POPCAT"
$ws.Range("C44").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D44").Value = "'2.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +50.67%  "
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").Value = "'1.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.53%  "
$ws.Range("D47").Value = "'138.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").Value = "'0.505"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.67%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'0.572"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "
